$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-05-14"

# Update the row label text for the May row
$ws.Range("A6").Value = "May (through 05-14)"

# Update May row (row 6) values for columns C..I (2016..2022)
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 29
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 17
$ws.Range("G6").Value = 28
$ws.Range("H6").Value = 53
$ws.Range("I6").Value = 53

# Update Total row (row 7) values for columns C..I (2016..2022)
$ws.Range("C7").Value = 182
$ws.Range("D7").Value = 282
$ws.Range("E7").Value = 265
$ws.Range("F7").Value = 172
$ws.Range("G7").Value = 290
$ws.Range("H7").Value = 576
$ws.Range("I7").Value = 605
